$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00009552326474482342
$ws.Range("C2").Value = 0.002658071450198252
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 0.685746420315646

$ws.Range("B3").Value = 0.6545652718822623
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 19.36876847130326
